# Auto-generated script to apply TPM data update to Comp-Cd47 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 1.254529666666667
$ws.Cells.Item(2, 8).Value2 = 3.763589
$ws.Cells.Item(2, 9).Value2 = 0.01270475613604028
$ws.Cells.Item(2, 10).Value2 = 0.01270475613604028
$ws.Cells.Item(2, 13).Value2 = 55.848606
$ws.Cells.Item(2, 14).Value2 = 167.545818
$ws.Cells.Item(2, 15).Value2 = 0.2323375192077237
$ws.Cells.Item(2, 16).Value2 = 0.2323375192077236
$ws.Cells.Item(2, 17).Value2 = 70.063733068978
$ws.Cells.Item(2, 18).Value2 = 630.573597620802
$ws.Cells.Item(2, 19).Value2 = 0.002951791522786703
$ws.Cells.Item(2, 20).Value2 = 0.002951791522786702
$ws.Cells.Item(3, 7).Value2 = 1.254529666666667
$ws.Cells.Item(3, 8).Value2 = 3.763589
$ws.Cells.Item(3, 9).Value2 = 0.01270475613604028
$ws.Cells.Item(3, 10).Value2 = 0.01270475613604028
$ws.Cells.Item(3, 15).Value2 = 0.3515710112922583
$ws.Cells.Item(3, 16).Value2 = 0.3515710112922583
$ws.Cells.Item(3, 17).Value2 = 106.0198007363099
$ws.Cells.Item(3, 18).Value2 = 954.178206626789
$ws.Cells.Item(3, 19).Value2 = 0.004466623962969204
$ws.Cells.Item(3, 20).Value2 = 0.004466623962969203
$ws.Cells.Item(4, 7).Value2 = 1.254529666666667
$ws.Cells.Item(4, 8).Value2 = 3.763589
$ws.Cells.Item(4, 9).Value2 = 0.01270475613604028
$ws.Cells.Item(4, 10).Value2 = 0.01270475613604028
$ws.Cells.Item(4, 13).Value2 = 33.195992
$ws.Cells.Item(4, 14).Value2 = 99.58797600000001
$ws.Cells.Item(4, 15).Value2 = 0.1380996766314891
$ws.Cells.Item(4, 16).Value2 = 0.1380996766314891
$ws.Cells.Item(4, 17).Value2 = 41.64535677842934
$ws.Cells.Item(4, 18).Value2 = 374.8082110058641
$ws.Cells.Item(4, 19).Value2 = 0.001754522714069089
$ws.Cells.Item(4, 20).Value2 = 0.001754522714069089
$ws.Cells.Item(5, 7).Value2 = 1.254529666666667
$ws.Cells.Item(5, 8).Value2 = 3.763589
$ws.Cells.Item(5, 9).Value2 = 0.01270475613604028
$ws.Cells.Item(5, 10).Value2 = 0.01270475613604028
$ws.Cells.Item(5, 13).Value2 = 66.82284533333335
$ws.Cells.Item(5, 14).Value2 = 200.468536
$ws.Cells.Item(5, 15).Value2 = 0.277991792868529
$ws.Cells.Item(5, 16).Value2 = 0.2779917928685289
$ws.Cells.Item(5, 17).Value2 = 83.83124188174492
$ws.Cells.Item(5, 18).Value2 = 754.4811769357041
$ws.Cells.Item(5, 19).Value2 = 0.003531817936215282
$ws.Cells.Item(5, 20).Value2 = 0.00353181793621528
$ws.Cells.Item(6, 9).Value2 = 0.9734152842234517
$ws.Cells.Item(6, 10).Value2 = 0.9734152842234516
$ws.Cells.Item(6, 13).Value2 = 55.848606
$ws.Cells.Item(6, 14).Value2 = 167.545818
$ws.Cells.Item(6, 15).Value2 = 0.2323375192077237
$ws.Cells.Item(6, 16).Value2 = 0.2323375192077236
$ws.Cells.Item(6, 17).Value2 = 5368.155666178074
$ws.Cells.Item(6, 18).Value2 = 48313.40099560267
$ws.Cells.Item(6, 19).Value2 = 0.226160892295358
$ws.Cells.Item(6, 20).Value2 = 0.2261608922953579
$ws.Cells.Item(7, 9).Value2 = 0.9734152842234517
$ws.Cells.Item(7, 10).Value2 = 0.9734152842234516
$ws.Cells.Item(7, 15).Value2 = 0.3515710112922583
$ws.Cells.Item(7, 16).Value2 = 0.3515710112922583
$ws.Cells.Item(7, 19).Value2 = 0.3422245958817799
$ws.Cells.Item(7, 20).Value2 = 0.3422245958817799
$ws.Cells.Item(8, 9).Value2 = 0.9734152842234517
$ws.Cells.Item(8, 10).Value2 = 0.9734152842234516
$ws.Cells.Item(8, 13).Value2 = 33.195992
$ws.Cells.Item(8, 14).Value2 = 99.58797600000001
$ws.Cells.Item(8, 15).Value2 = 0.1380996766314891
$ws.Cells.Item(8, 16).Value2 = 0.1380996766314891
$ws.Cells.Item(8, 17).Value2 = 3190.791414725769
$ws.Cells.Item(8, 18).Value2 = 28717.12273253192
$ws.Cells.Item(8, 19).Value2 = 0.1344283359794078
$ws.Cells.Item(8, 20).Value2 = 0.1344283359794077
$ws.Cells.Item(9, 9).Value2 = 0.9734152842234517
$ws.Cells.Item(9, 10).Value2 = 0.9734152842234516
$ws.Cells.Item(9, 13).Value2 = 66.82284533333335
$ws.Cells.Item(9, 14).Value2 = 200.468536
$ws.Cells.Item(9, 15).Value2 = 0.277991792868529
$ws.Cells.Item(9, 16).Value2 = 0.2779917928685289
$ws.Cells.Item(9, 17).Value2 = 6422.997125591183
$ws.Cells.Item(9, 18).Value2 = 57806.97413032065
$ws.Cells.Item(9, 19).Value2 = 0.2706014600669061
$ws.Cells.Item(9, 20).Value2 = 0.270601460066906
$ws.Cells.Item(10, 7).Value2 = 1.151276666666667
$ws.Cells.Item(10, 8).Value2 = 3.45383
$ws.Cells.Item(10, 9).Value2 = 0.01165910195968263
$ws.Cells.Item(10, 10).Value2 = 0.01165910195968263
$ws.Cells.Item(10, 13).Value2 = 55.848606
$ws.Cells.Item(10, 14).Value2 = 167.545818
$ws.Cells.Item(10, 15).Value2 = 0.2323375192077237
$ws.Cells.Item(10, 16).Value2 = 0.2323375192077236
$ws.Cells.Item(10, 17).Value2 = 64.29719695366001
$ws.Cells.Item(10, 18).Value2 = 578.67477258294
$ws.Cells.Item(10, 19).Value2 = 0.002708846825502572
$ws.Cells.Item(10, 20).Value2 = 0.002708846825502571
$ws.Cells.Item(11, 7).Value2 = 1.151276666666667
$ws.Cells.Item(11, 8).Value2 = 3.45383
$ws.Cells.Item(11, 9).Value2 = 0.01165910195968263
$ws.Cells.Item(11, 10).Value2 = 0.01165910195968263
$ws.Cells.Item(11, 15).Value2 = 0.3515710112922583
$ws.Cells.Item(11, 16).Value2 = 0.3515710112922583
$ws.Cells.Item(11, 17).Value2 = 97.29393097309222
$ws.Cells.Item(11, 18).Value2 = 875.6453787578299
$ws.Cells.Item(11, 19).Value2 = 0.004099002266725172
$ws.Cells.Item(11, 20).Value2 = 0.004099002266725171
$ws.Cells.Item(12, 7).Value2 = 1.151276666666667
$ws.Cells.Item(12, 8).Value2 = 3.45383
$ws.Cells.Item(12, 9).Value2 = 0.01165910195968263
$ws.Cells.Item(12, 10).Value2 = 0.01165910195968263
$ws.Cells.Item(12, 13).Value2 = 33.195992
$ws.Cells.Item(12, 14).Value2 = 99.58797600000001
$ws.Cells.Item(12, 15).Value2 = 0.1380996766314891
$ws.Cells.Item(12, 16).Value2 = 0.1380996766314891
$ws.Cells.Item(12, 17).Value2 = 38.21777101645334
$ws.Cells.Item(12, 18).Value2 = 343.95993914808
$ws.Cells.Item(12, 19).Value2 = 0.001610118210445732
$ws.Cells.Item(12, 20).Value2 = 0.001610118210445732
$ws.Cells.Item(13, 7).Value2 = 1.151276666666667
$ws.Cells.Item(13, 8).Value2 = 3.45383
$ws.Cells.Item(13, 9).Value2 = 0.01165910195968263
$ws.Cells.Item(13, 10).Value2 = 0.01165910195968263
$ws.Cells.Item(13, 13).Value2 = 66.82284533333335
$ws.Cells.Item(13, 14).Value2 = 200.468536
$ws.Cells.Item(13, 15).Value2 = 0.277991792868529
$ws.Cells.Item(13, 16).Value2 = 0.2779917928685289
$ws.Cells.Item(13, 17).Value2 = 76.93158263254224
$ws.Cells.Item(13, 18).Value2 = 692.3842436928801
$ws.Cells.Item(13, 19).Value2 = 0.003241134657009154
$ws.Cells.Item(13, 20).Value2 = 0.003241134657009153
$ws.Cells.Item(14, 7).Value2 = 0.2192983333333333
$ws.Cells.Item(14, 8).Value2 = 0.6578949999999999
$ws.Cells.Item(14, 9).Value2 = 0.002220857680825461
$ws.Cells.Item(14, 10).Value2 = 0.002220857680825461
$ws.Cells.Item(14, 13).Value2 = 55.848606
$ws.Cells.Item(14, 14).Value2 = 167.545818
$ws.Cells.Item(14, 15).Value2 = 0.2323375192077237
$ws.Cells.Item(14, 16).Value2 = 0.2323375192077236
$ws.Cells.Item(14, 17).Value2 = 12.24750621479
$ws.Cells.Item(14, 18).Value2 = 110.22755593311
$ws.Cells.Item(14, 19).Value2 = 0.0005159885640764062
$ws.Cells.Item(14, 20).Value2 = 0.000515988564076406
$ws.Cells.Item(15, 7).Value2 = 0.2192983333333333
$ws.Cells.Item(15, 8).Value2 = 0.6578949999999999
$ws.Cells.Item(15, 9).Value2 = 0.002220857680825461
$ws.Cells.Item(15, 10).Value2 = 0.002220857680825461
$ws.Cells.Item(15, 15).Value2 = 0.3515710112922583
$ws.Cells.Item(15, 16).Value2 = 0.3515710112922583
$ws.Cells.Item(15, 17).Value2 = 18.5328145037661
$ws.Cells.Item(15, 18).Value2 = 166.795330533895
$ws.Cells.Item(15, 19).Value2 = 0.0007807891807839867
$ws.Cells.Item(15, 20).Value2 = 0.0007807891807839866
$ws.Cells.Item(16, 7).Value2 = 0.2192983333333333
$ws.Cells.Item(16, 8).Value2 = 0.6578949999999999
$ws.Cells.Item(16, 9).Value2 = 0.002220857680825461
$ws.Cells.Item(16, 10).Value2 = 0.002220857680825461
$ws.Cells.Item(16, 13).Value2 = 33.195992
$ws.Cells.Item(16, 14).Value2 = 99.58797600000001
$ws.Cells.Item(16, 15).Value2 = 0.1380996766314891
$ws.Cells.Item(16, 16).Value2 = 0.1380996766314891
$ws.Cells.Item(16, 17).Value2 = 7.279825718946666
$ws.Cells.Item(16, 18).Value2 = 65.51843147052
$ws.Cells.Item(16, 19).Value2 = 0.0003066997275665551
$ws.Cells.Item(16, 20).Value2 = 0.0003066997275665549
$ws.Cells.Item(17, 7).Value2 = 0.2192983333333333
$ws.Cells.Item(17, 8).Value2 = 0.6578949999999999
$ws.Cells.Item(17, 9).Value2 = 0.002220857680825461
$ws.Cells.Item(17, 10).Value2 = 0.002220857680825461
$ws.Cells.Item(17, 13).Value2 = 66.82284533333335
$ws.Cells.Item(17, 14).Value2 = 200.468536
$ws.Cells.Item(17, 15).Value2 = 0.277991792868529
$ws.Cells.Item(17, 16).Value2 = 0.2779917928685289
$ws.Cells.Item(17, 17).Value2 = 14.65413861019111
$ws.Cells.Item(17, 18).Value2 = 131.88724749172
$ws.Cells.Item(17, 19).Value2 = 0.0006173802083985132
$ws.Cells.Item(17, 20).Value2 = 0.000617380208398513
